$d = $word.ActiveDocument

# --- Helper: rebuild a paragraph's runs as a single merged run with en-US language ---
function MergeParagraphRuns($paraIndex, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $textRng = $d.Range($rng.Start, $rng.End - 1)
    $textRng.Delete()
    $ins = $d.Range($rng.Start, $rng.Start)
    $ins.InsertAfter($newText)
    $ins.LanguageID = "en-US"
}

# 1. Add language tagging to the title (Heading1) run.
$d.Paragraphs(1).Range.LanguageID = "en-US"

# 2. "GET /gateslots/reservations/{reservationNo}" paragraph: collapse its
#    3 runs ("GET ", "/gateslots/reservations", "/{reservationNo}") into one run.
MergeParagraphRuns 10 "GET /gateslots/reservations/{reservationNo}"

# 3. "GET /gateslots/reservations?warehouseName=...&loadDate=..." paragraph: same collapse.
MergeParagraphRuns 11 "GET /gateslots/reservations?warehouseName=…&loadDate=…"

# 4. Move the "_GoBack" bookmark away from the end of the POST /login paragraph -
#    it will be re-added further down, at the new edit location.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 5. Turn the trailing empty paragraph into the new
#    "GET /gateslots/reservations/sessions/{sessionId}" entry, split into
#    runs the same way the real edit left them, with the _GoBack bookmark
#    sitting right before the closing "}".
$p12 = $d.Paragraphs(12)
$start = $p12.Range.Start

$part1 = "GET /gateslots/reservations"
$r1 = $d.Range($start, $start)
$r1.InsertAfter($part1)
$r1.LanguageID = "en-US"

$afterPart1 = $start + $part1.Length
$part2 = "/sessions/{sessionId"
$r2 = $d.Range($afterPart1, $afterPart1)
$r2.InsertAfter($part2)
# Toggling Bold on/off keeps this run distinct from the previous one instead
# of silently coalescing them on save.
$r2.Bold = 1
$r2.LanguageID = "en-US"
$r2.Bold = 0

$afterPart2 = $afterPart1 + $part2.Length
$part3 = "}"
$r3 = $d.Range($afterPart2, $afterPart2)
$r3.InsertAfter($part3)
$r3.LanguageID = "en-US"

# Bookmark goes between "/sessions/{sessionId" and "}".
$bmRng = $d.Range($afterPart2, $afterPart2)
$d.Bookmarks.Add("_GoBack", $bmRng)
